$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture the existing data rows (2-21, columns A-H) before touching anything ---
$oldData = @{}
for ($r = 2; $r -le 21; $r++) {
    $rowVals = @()
    for ($c = 1; $c -le 8; $c++) {
        $rowVals += ,$ws.Cells.Item($r, $c).Value()
    }
    $oldData[$r] = $rowVals
}

# --- Shift the old rows down by 6 (old row 2 -> new row 8, ... old row 21 -> new row 27) ---
# Column A (timestamp) always follows (row-2)*100 positionally, so it is
# recomputed for the new row rather than carried verbatim from the old row.
for ($r = 21; $r -ge 2; $r--) {
    $target = $r + 6
    $vals = $oldData[$r]
    $ws.Cells.Item($target, 1).Value = ($target - 2) * 100
    for ($c = 2; $c -le 8; $c++) {
        $ws.Cells.Item($target, $c).Value = $vals[$c - 1]
    }
}

# --- New rows inserted at the top (rows 2-7) ---
$newTop = @{
    2  = @(0,   -0.3537254333496094, 0.2088937759399415, 0.5035260319709778, -0.0386372283101081, 0.008399397134780801, -0.0021380283869802)
    3  = @(100, -0.4043011069297791, 0.3229363560676575, 0.4749223440885544, 0.0134390350431203, 0.0704022198915481, -0.0390953756868839)
    4  = @(200, -0.2195036411285398, 0.2696369886398314, 0.4450621306896209, 0.0154243474826216, 0.030695978552103, 0.0099265603348612)
    5  = @(300, -0.5155707597732552, 0.2643678188323975, 0.5658968165516856, 0.041233405470848, -0.0022907445672899, 0.0502436682581901)
    6  = @(400, -0.4721715450286855, 0.2206716537475583, 0.4629700779914848, -0.030695978552103, -0.062460970133543, 0.0204639863222837)
    7  = @(500, -0.2396689057350159, 0.2092438936233521, 0.5360905304551128, 0.0216857157647609, -0.0343611687421798, -0.0035124751739203)
}

# --- New rows appended at the bottom (rows 28-31) ---
$newBottom = @{
    28 = @(2600, -0.1972274780273441, 0.3113194406032568, 0.5224930047988898, 0.0047342055477201, 0.0751364231109619, 0.0545197241008281)
    29 = @(2700, -0.2138409614562988, 0.2108606994152066, 0.46574055776, 0.0003054326225537, 0.0100792767480015, -0.0006108652451075)
    30 = @(2800, -0.09080266952514603, 0.2757070064544683, 0.4601370841264726, -0.0152716310694813, -0.00167987938039, 0.047036625444889)
    31 = @(2900, -0.04834830760955861, 0.3611972928047176, 0.5197352617979051, -0.0091629782691597, -0.0114537235349416, 0.0062613687478005)
}

function Write-DataRow($rowNum, $vals) {
    $ws.Cells.Item($rowNum, 1).Value = $vals[0]
    $ws.Cells.Item($rowNum, 2).Value = "falling"
    $ws.Cells.Item($rowNum, 3).Value = $vals[1]
    $ws.Cells.Item($rowNum, 4).Value = $vals[2]
    $ws.Cells.Item($rowNum, 5).Value = $vals[3]
    $ws.Cells.Item($rowNum, 6).Value = $vals[4]
    $ws.Cells.Item($rowNum, 7).Value = $vals[5]
    $ws.Cells.Item($rowNum, 8).Value = $vals[6]
}

foreach ($r in ($newTop.Keys | Sort-Object)) {
    Write-DataRow $r $newTop[$r]
}

foreach ($r in ($newBottom.Keys | Sort-Object)) {
    Write-DataRow $r $newBottom[$r]
}

Write-Host "Final UsedRange:" $ws.UsedRange.Address()
